$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

# Row 2
Set-TextValue "D2" "68.235.08"
$ws.Range("E2").Value = "  +2.77%  "

# Row 3
Set-TextValue "D3" "3.631.39"
$ws.Range("E3").Value = "  +2.04%  "

# Row 4
$ws.Range("E4").Value = "  +0.12%  "

# Row 5
Set-TextValue "D5" "197.20"
$ws.Range("E5").Value = "  +9.22%  "

# Row 6
Set-TextValue "D6" "575.14"
$ws.Range("E6").Value = "  -1.52%  "

# Row 7
Set-TextValue "D7" "3.626.12"
$ws.Range("E7").Value = "  +2.19%  "

# Row 8
$ws.Range("E8").Value = "  +2.34%  "

# Row 9
$ws.Range("E9").Value = "  -0.17%  "

# Row 10
Set-TextValue "D10" "0.676"
$ws.Range("E10").Value = "  +1.50%  "

# Row 11
$ws.Range("E11").Value = "  +8.15%  "

# Row 12
Set-TextValue "D12" "56.51"
$ws.Range("E12").Value = "  +5.90%  "

# Row 13
Set-TextValue "D13" "0.0000295"
$ws.Range("E13").Value = "  +17.66%  "

# Row 14
$ws.Range("E14").Value = "  +2.51%  "

# Row 15
Set-TextValue "D15" "4.213.63"
$ws.Range("E15").Value = "  +1.84%  "

# Row 16
Set-TextValue "D16" "3.632.84"
$ws.Range("E16").Value = "  +2.03%  "

# Row 17
$ws.Range("E17").Value = "  +0.71%  "

# Row 18
Set-TextValue "D18" "12.54"
$ws.Range("E18").Value = "  +4.08%  "

# Row 19
Set-TextValue "D19" "68.149.19"
$ws.Range("E19").Value = "  +3.00%  "

# Row 20
Set-TextValue "D20" "18.57"
$ws.Range("E20").Value = "  +2.26%  "

# Row 21
$ws.Range("E21").Value = "  +3.77%  "

# Row 22
Set-TextValue "D22" "403.01"
$ws.Range("E22").Value = "  +3.32%  "

# Row 23
Set-TextValue "D23" "13.30"
$ws.Range("E23").Value = "  +30.97%  "

# Row 24
Set-TextValue "D24" "4.24"
$ws.Range("E24").Value = "  -0.36%  "

# Row 25
Set-TextValue "D25" "85.93"
$ws.Range("E25").Value = "  +1.95%  "

# Row 26
Set-TextValue "D26" "2.97"
$ws.Range("E26").Value = "  +4.20%  "

# Row 27
Set-TextValue "D27" "12.65"
$ws.Range("E27").Value = "  +4.34%  "

# Row 28
Set-TextValue "D28" "3.88"
$ws.Range("E28").Value = "  +7.96%  "

# Row 29
$ws.Range("E29").Value = "  +2.28%  "

# Row 30
Set-TextValue "D30" "8.24"
$ws.Range("E30").Value = "  +23.45%  "

# Row 31
Set-TextValue "D31" "9.21"
$ws.Range("E31").Value = "  +3.75%  "

# Row 32
Set-TextValue "D32" "31.84"
$ws.Range("E32").Value = "  +3.07%  "

# Row 33
Set-TextValue "D33" "690.66"
$ws.Range("E33").Value = "  +16.32%  "

# Row 34
$ws.Range("E34").Value = "  +3.26%  "

# Row 35
$ws.Range("E35").Value = "  +5.21%  "

# Row 36
Set-TextValue "D36" "64.75"
$ws.Range("E36").Value = "  -0.65%  "

# Row 37
Set-TextValue "D37" "42.77"
$ws.Range("E37").Value = "  +3.64%  "

# Row 38
Set-TextValue "D38" "0.427"
$ws.Range("E38").Value = "  +15.92%  "

# Row 39
$ws.Range("E39").Value = "  +0.04%  "

# Row 40
Set-TextValue "D40" "0.0₃0790"
$ws.Range("E40").Value = "  +8.19%  "

# Row 41 (was Kaspa -> now Fetch.AI)
$ws.Range("B41").Value = "Fetch.AI"
$ws.Range("C41").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D41" "2.92"
$ws.Range("E41").Value = "  +22.40%  "

# Row 42 (was Fetch.AI -> now Kaspa)
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D42" "0.139"
$ws.Range("E42").Value = "  +8.21%  "

# Row 43
Set-TextValue "D43" "3.16"
$ws.Range("E43").Value = "  +15.47%  "

# Row 44
Set-TextValue "D44" "3.231.11"
$ws.Range("E44").Value = "  +16.21%  "

# Row 45
Set-TextValue "D45" "3.12"
$ws.Range("E45").Value = "  +43.75%  "

# Row 46
Set-TextValue "D46" "0.998"
$ws.Range("E46").Value = "  -0.09%  "

# Row 47
$ws.Range("E47").Value = "  +3.65%  "

# Row 48
Set-TextValue "D48" "8.95"
$ws.Range("E48").Value = "  +9.71%  "

# Row 49
$ws.Range("E49").Value = "  +2.14%  "

# Row 50
Set-TextValue "D50" "3.10"
$ws.Range("E50").Value = "  +1.30%  "

# Row 51 (was WEMIXToken -> now Monero)
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D51" "142.52"
$ws.Range("E51").Value = "  +5.72%  "
